$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value2 = -4025.31637494483
$ws.Range("D3").Value2 = -6794.77939906494
$ws.Range("E3").Value2 = 6437.96399817819
$ws.Range("F3").Value2 = 9207.4270222983
$ws.Range("C4").Value2 = -300452.377158042
$ws.Range("D4").Value2 = -460312.812688958
$ws.Range("E4").Value2 = 303514.72511477
$ws.Range("F4").Value2 = 463375.160645686
$ws.Range("C5").Value2 = -11207563.8950994
$ws.Range("D5").Value2 = -17141236.76244
$ws.Range("E5").Value2 = 11210385.834741
$ws.Range("F5").Value2 = 17144058.7020816
$ws.Range("C6").Value2 = -468477934.473888
$ws.Range("D6").Value2 = -716475924.489544
$ws.Range("E6").Value2 = 468480775.176764
$ws.Range("F6").Value2 = 716478765.19242
$ws.Range("C7").Value2 = -18437084807.2648
$ws.Range("D7").Value2 = -28197088698.3042
$ws.Range("E7").Value2 = 18437087646.1339
$ws.Range("F7").Value2 = 28197091537.1732
$ws.Range("C8").Value2 = -769709026033.162
$ws.Range("D8").Value2 = -1177168372456.67
$ws.Range("E8").Value2 = 769709029072.728
$ws.Range("F8").Value2 = 1177168375496.24
$ws.Range("C9").Value2 = -31427073643579.5
$ws.Range("D9").Value2 = -48063561496880.6
$ws.Range("E9").Value2 = 31427073646742.0
$ws.Range("F9").Value2 = 48063561500043.2
$ws.Range("C10").Value2 = -1172202265538391.0
$ws.Range("D10").Value2 = -1792728661742270.0
$ws.Range("E10").Value2 = 1172202265541279.0
$ws.Range("F10").Value2 = 1792728661745158.0
$ws.Range("C11").Value2 = -46453136203764776.0
$ws.Range("D11").Value2 = -71043941091533832.0
$ws.Range("E11").Value2 = 46453136203767528.0
$ws.Range("F11").Value2 = 71043941091536584.0
$ws.Range("C12").Value2 = -1623970856810717184.0
$ws.Range("D12").Value2 = -2483649099159784448.0
$ws.Range("E12").Value2 = 1623970856810719744.0
$ws.Range("F12").Value2 = 2483649099159787520.0
$ws.Range("C13").Value2 = -50482943368827158528.0
$ws.Range("D13").Value2 = -77206999309800980480.0
$ws.Range("E13").Value2 = 50482943368827158528.0
$ws.Range("F13").Value2 = 77206999309800980480.0
$ws.Range("C14").Value2 = -1829844178826780475392.0
$ws.Range("D14").Value2 = -2798505174699459739648.0
$ws.Range("E14").Value2 = 1829844178826780475392.0
$ws.Range("F14").Value2 = 2798505174699459739648.0
$ws.Range("C15").Value2 = -68452498712492070928384.0
$ws.Range("D15").Value2 = -104689062645126666387456.0
$ws.Range("E15").Value2 = 68452498712492070928384.0
$ws.Range("F15").Value2 = 104689062645126666387456.0
$ws.Range("C16").Value2 = -3055876274041987096641536.0
$ws.Range("D16").Value2 = -4673559456647788248432640.0
$ws.Range("E16").Value2 = 3055876274041987096641536.0
$ws.Range("F16").Value2 = 4673559456647788248432640.0
$ws.Range("C17").Value2 = -110022277208659611141275648.0
$ws.Range("D17").Value2 = -168264552612377982503223296.0
$ws.Range("E17").Value2 = 110022277208659611141275648.0
$ws.Range("F17").Value2 = 168264552612377982503223296.0
$ws.Range("C18").Value2 = -4749869741741995097329762304.0
$ws.Range("D18").Value2 = -7264298897808871624452079616.0
$ws.Range("E18").Value2 = 4749869741741995097329762304.0
$ws.Range("F18").Value2 = 7264298897808871624452079616.0
$ws.Range("C19").Value2 = -133581122853225166942830067712.0
$ws.Range("D19").Value2 = -204294697806781230722991521792.0
$ws.Range("E19").Value2 = 133581122853225166942830067712.0
$ws.Range("F19").Value2 = 204294697806781230722991521792.0
$ws.Range("C20").Value2 = -4227640411149868047962954792960.0
$ws.Range("D20").Value2 = -6465618058777577272324747427840.0
$ws.Range("E20").Value2 = 4227640411149868047962954792960.0
$ws.Range("F20").Value2 = 6465618058777577272324747427840.0
$ws.Range("C21").Value2 = -108318139305260009422033617158144.0
$ws.Range("D21").Value2 = -165658298595643636334353721589760.0
$ws.Range("E21").Value2 = 108318139305260009422033617158144.0
$ws.Range("F21").Value2 = 165658298595643636334353721589760.0
$ws.Range("C22").Value2 = -3645650571580160727307625493954560.0
$ws.Range("D22").Value2 = -5575541408260488085701211275132928.0
$ws.Range("E22").Value2 = 3645650571580160727307625493954560.0
$ws.Range("F22").Value2 = 5575541408260488085701211275132928.0
$ws.Range("C23").Value2 = -96657382565018871946263428211933184.0
$ws.Range("D23").Value2 = -147824710109765476273038241470349312.0
$ws.Range("E23").Value2 = 96657382565018871946263428211933184.0
$ws.Range("F23").Value2 = 147824710109765476273038241470349312.0
$ws.Range("C24").Value2 = -2431490268428494140142706216816082944.0
$ws.Range("D24").Value2 = -3718643465473275983272872739052978176.0
$ws.Range("E24").Value2 = 2431490268428494140142706216816082944.0
$ws.Range("F24").Value2 = 3718643465473275983272872739052978176.0
$ws.Range("C25").Value2 = -79708875951275134778005129961406464000.0
$ws.Range("D25").Value2 = -121904206052201243438846507934510219264.0
$ws.Range("E25").Value2 = 79708875951275134778005129961406464000.0
$ws.Range("F25").Value2 = 121904206052201243438846507934510219264.0
$ws.Range("C26").Value2 = -3330818176562237738227308012185581518848.0
$ws.Range("D26").Value2 = -5094046810624540044991039348742924795904.0
$ws.Range("E26").Value2 = 3330818176562237738227308012185581518848.0
$ws.Range("F26").Value2 = 5094046810624540044991039348742924795904.0
$ws.Range("C27").Value2 = -74439569649861466866813883362127407742976.0
$ws.Range("D27").Value2 = -113845497489903181179227214722871945330688.0
$ws.Range("E27").Value2 = 74439569649861466866813883362127407742976.0
$ws.Range("F27").Value2 = 113845497489903181179227214722871945330688.0
$ws.Range("C28").Value2 = -2990391831475307390890294864089891232284672.0
$ws.Range("D28").Value2 = -4573409644163393687295274513576526945976320.0
$ws.Range("E28").Value2 = 2990391831475307390890294864089891232284672.0
$ws.Range("F28").Value2 = 4573409644163393687295274513576526945976320.0
$ws.Range("C29").Value2 = -68565725390996450355809020622540903739817984.0
$ws.Range("D29").Value2 = -104862227906614510755686366394778939593588736.0
$ws.Range("E29").Value2 = 68565725390996450355809020622540903739817984.0
$ws.Range("F29").Value2 = 104862227906614510755686366394778939593588736.0
$ws.Range("C30").Value2 = -1885528820138701285121106619809625900125782016.0
$ws.Range("D30").Value2 = -2883667484510237508906965786354625650458886144.0
$ws.Range("E30").Value2 = 1885528820138701285121106619809625900125782016.0
$ws.Range("F30").Value2 = 2883667484510237508906965786354625650458886144.0
$ws.Range("C31").Value2 = -38589770227748874450232071572329223627187683328.0
$ws.Range("D31").Value2 = -59017960612394548787198923429376196160031555584.0
$ws.Range("E31").Value2 = 38589770227748874450232071572329223627187683328.0
$ws.Range("F31").Value2 = 59017960612394548787198923429376196160031555584.0
$ws.Range("C32").Value2 = -936136256593958261651706639792303269239425335296.0
$ws.Range("D32").Value2 = -1431696856276400672932393088763144260558910390272.0
$ws.Range("E32").Value2 = 936136256593958261651706639792303269239425335296.0
$ws.Range("F32").Value2 = 1431696856276400672932393088763144260558910390272.0
$ws.Range("C33").Value2 = -18790439503736670252152112373949584417178869825536.0
$ws.Range("D33").Value2 = -28737497320563992804495032987669498439993539952640.0
$ws.Range("E33").Value2 = 18790439503736670252152112373949584417178869825536.0
$ws.Range("F33").Value2 = 28737497320563992804495032987669498439993539952640.0
$ws.Range("C34").Value2 = -574001933920375720303231458852719676412934051856384.0
$ws.Range("D34").Value2 = -877860203043950852768080585487115716396610291236864.0
$ws.Range("E34").Value2 = 574001933920375720303231458852719676412934051856384.0
$ws.Range("F34").Value2 = 877860203043950852768080585487115716396610291236864.0
$ws.Range("C35").Value2 = -9963397492177624209980509222908067790725433645858816.0
$ws.Range("D35").Value2 = -15237701527855706204362549319075185370046840946294784.0
$ws.Range("E35").Value2 = 9963397492177624209980509222908067790725433645858816.0
$ws.Range("F35").Value2 = 15237701527855706204362549319075185370046840946294784.0
$ws.Range("C36").Value2 = -221393539010899506685292182370271968800778054065455104.0
$ws.Range("D36").Value2 = -338592199126087347304523770793835838262860908991610880.0
$ws.Range("E36").Value2 = 221393539010899506685292182370271968800778054065455104.0
$ws.Range("F36").Value2 = 338592199126087347304523770793835838262860908991610880.0

Write-Host "Done updating cells"
